# Update the price list worksheet: bump the date in A1 and the four
# package prices in column D (rows 29-32).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45436

$ws.Range("D29").Value = 1067
$ws.Range("D30").Value = 1265
$ws.Range("D31").Value = 2451
$ws.Range("D32").Value = 2750
